$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: relabel Model b$/PJ -> Model bGC/PJ, $/GJ -> GC/GJ, add new $/GJ column ---
$ws.Range("O4").Value = "Model bGC/PJ"
$ws.Range("P4").Value = "GC/GJ"
$ws.Range("Q4").Value = "$/GJ"

# --- New row 2: exchange rate ---
$ws.Range("I2").Value = "exchange rate"
$ws.Range("J2").Formula = "=0.541/0.147"

# --- Row 18: add GC/GJ conversions for 216 and 130 ---
$ws.Range("Q18").Formula = "=216/3.7"
$ws.Range("Q18").NumberFormat = "0.0"
$ws.Range("R18").Formula = "=130/3.7"
$ws.Range("R18").NumberFormat = "0"

# --- Row 27: relabel Model b$/PJ -> Model bGC/PJ, $/GJ -> GC/GJ, add new $/GJ column ---
$ws.Range("O27").Value = "Model bGC/PJ"
$ws.Range("P27").Value = "GC/GJ"
$ws.Range("Q27").Value = "$/GJ"

# --- Row 28: add GC/GJ conversion ---
$ws.Range("Q28").Formula = "=P28/3.7"

# --- Row 50: relabel Model b$/PJ -> Model bGC/PJ, $/GJ -> GC/GJ, add new $/GJ column ---
$ws.Range("O50").Value = "Model bGC/PJ"
$ws.Range("P50").Value = "GC/GJ"
$ws.Range("Q50").Value = "$/GJ"

# --- Row 51: fix formulas (divide instead of multiply, /100 correction) and add GC/GJ conversion ---
$ws.Range("M51").Formula = "=0.147/0.541*L51/100"
$ws.Range("N51").Formula = "=M51/N49"
$ws.Range("Q51").Formula = "=P51/3.7"

# --- View tweaks ---
$ws.Range("N51").Select()
